$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: add new prefix mapping "Item starts with AMP use:" -> shared template file
$ws.Range("A9").Value = "Item starts with AMP use:"

# Row 10: add new prefix mapping "Item starts with ASP use:" -> same shared template file
# Match the style used by the rest of the data rows (style index 2 / vertical-top wrap)
# instead of the bold/centered header style A10 previously had
$ws.Range("B10").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A10").Value = "Item starts with ASP use:"

$ws.Range("B9").Value = "F-825-1100A CMD Final_In Process Inspection Template_final.xlsx"
$ws.Range("B10").Value = "F-825-1100A CMD Final_In Process Inspection Template_final.xlsx"

# Row 9 grew taller to fit the wrapped text
$ws.Rows.Item(9).RowHeight = 30

# Reflect the last-used selection in the saved view state
[void]$ws.Range("C10").Select()
